$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-21, columns A (1), B (2), C (3)
$data = @(
    @("航天电子", "利欧股份", "利欧股份"),
    @("中国卫星", "航天电子", "锋龙股份"),
    @("利欧股份", "金风科技", "航天电子"),
    @("国晟科技", "白银有色", "金风科技"),
    @("金风科技", "湖南白银", "白银有色"),
    @("乾照光电", "岩山科技", "国晟科技"),
    @("巨力索具", "国晟科技", "隆基绿能"),
    @("白银有色", "锋龙股份", "中国卫星"),
    @("锋龙股份", "隆基绿能", "航天发展"),
    @("湖南白银", "蓝色光标", "蓝色光标"),
    @("信维通信", "中国卫星", "通宇通讯"),
    @("东方日升", "乾照光电", "湖南白银"),
    @("岩山科技", "巨力索具", "乾照光电"),
    @("隆基绿能", "航天发展", "巨力索具"),
    @("蓝色光标", "中国长城", "中国长城"),
    @("中国卫通", "特变电工", "岩山科技"),
    @("航天发展", "中国西电", "东方日升"),
    @("中国长城", "中国卫通", "博士眼镜"),
    @("钧达股份", "信维通信", "中国卫通"),
    @("特变电工", "平潭发展", "钧达股份")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
